# Applies the numeric updates described by the commit diff for Sheets/Exodus_Profits.xlsx
# (columns H:N = price/profit metrics recomputed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 334.875
$ws.Range("I39").Value = 239.8
$ws.Range("K39").Value = 719.4000000000001
$ws.Range("M39").Value = -423.4000000000001

$ws.Range("H43").Value = 3607
$ws.Range("I43").Value = 2942.8333
$ws.Range("J43").Value = 5599.5
$ws.Range("K43").Value = 2942.8333
$ws.Range("L43").Value = 5599.5
$ws.Range("M43").Value = -2873.8333
$ws.Range("N43").Value = -5737.5

$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").Value = ""

$ws.Range("H95").Value = 11124.8
$ws.Range("J95").Value = 11124.8
$ws.Range("L95").Value = 11124.8
$ws.Range("N95").Value = -16616.8

$ws.Range("H98").Value = 4237.1665
$ws.Range("I98").Value = 4237.1665
$ws.Range("K98").Value = 4237.1665
$ws.Range("M98").Value = -2739.1665

$ws.Range("H115").Value = 920
$ws.Range("I115").Value = 920
$ws.Range("K115").Value = 2760
$ws.Range("M115").Value = -1193

$ws.Range("H116").Value = 1860494
$ws.Range("I116").Value = 8650.666999999999
$ws.Range("K116").Value = 8650.666999999999
$ws.Range("M116").Value = -5208.666999999999

$ws.Range("H122").Value = 4237.1665
$ws.Range("I122").Value = 4237.1665
$ws.Range("K122").Value = 12711.4995
$ws.Range("M122").Value = -10261.4995

$ws.Range("H132").Value = 1290.5869
$ws.Range("I132").Value = 1290.5869
$ws.Range("K132").Value = 3871.7607
$ws.Range("M132").Value = -1341.7607

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2675.1875
$ws.Range("I63").Value = 2805
$ws.Range("J63").Value = 2285.75
$ws.Range("K63").Value = 2805
$ws.Range("L63").Value = 2285.75
$ws.Range("M63").Value = -2119
$ws.Range("N63").Value = -3657.75

$ws.Range("H66").Value = 2675.1875
$ws.Range("I66").Value = 2805
$ws.Range("J66").Value = 2285.75
$ws.Range("K66").Value = 14025
$ws.Range("L66").Value = 11428.75
$ws.Range("M66").Value = -10593
$ws.Range("N66").Value = -18292.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2565.3333
$ws.Range("I31").Value = 1848.4
$ws.Range("K31").Value = 1848.4
$ws.Range("M31").Value = -1553.4

$ws.Range("H34").Value = 2565.3333
$ws.Range("I34").Value = 1848.4
$ws.Range("K34").Value = 1848.4
$ws.Range("M34").Value = -1646.4

$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = ""

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").Value = ""

$ws.Range("H132").Value = 4002522.8
$ws.Range("I132").Value = 4764265
$ws.Range("J132").Value = 3375
$ws.Range("K132").Value = 14292795
$ws.Range("L132").Value = 10125
$ws.Range("M132").Value = -14290265
$ws.Range("N132").Value = -15185

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 660
$ws.Range("I113").Value = 408.57144
$ws.Range("K113").Value = 1225.71432
$ws.Range("M113").Value = 944.28568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 59996
$ws.Range("J15").Value = 59996
$ws.Range("L15").Value = 59996
$ws.Range("N15").Value = -60572

$ws.Range("H81").Value = 59996
$ws.Range("J81").Value = 59996
$ws.Range("L81").Value = 59996
$ws.Range("N81").Value = -61992

$ws.Range("H84").Value = 59996
$ws.Range("J84").Value = 59996
$ws.Range("L84").Value = 179988
$ws.Range("N84").Value = -189972

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2070.8
$ws.Range("I61").Value = 2233.8333
$ws.Range("K61").Value = 2233.8333
$ws.Range("M61").Value = -2031.8333

$ws.Range("H113").Value = 2070.8
$ws.Range("I113").Value = 2233.8333
$ws.Range("K113").Value = 2233.8333
$ws.Range("M113").Value = -63.83329999999978

$ws.Range("H132").Value = 12531.786
$ws.Range("J132").Value = 3715.8333
$ws.Range("L132").Value = 11147.4999
$ws.Range("N132").Value = -16207.4999

$ws.Range("H135").Value = 61999.2
$ws.Range("J135").Value = 61999.2
$ws.Range("L135").Value = 61999.2
$ws.Range("N135").Value = -72139.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 62422
$ws.Range("J27").Value = 62422
$ws.Range("L27").Value = 62422
$ws.Range("N27").Value = -62560

$ws.Range("H29").Value = 142500
$ws.Range("I29").Value = 275000
$ws.Range("J29").Value = 10000
$ws.Range("K29").Value = 275000
$ws.Range("L29").Value = 10000
$ws.Range("M29").Value = -274710
$ws.Range("N29").Value = -10580

$ws.Range("H81").Value = 43372.727
$ws.Range("I81").Value = 3567
$ws.Range("J81").Value = 58299.875
$ws.Range("K81").Value = 7134
$ws.Range("L81").Value = 116599.75
$ws.Range("M81").Value = -6073
$ws.Range("N81").Value = -118721.75

$ws.Range("H84").Value = 43372.727
$ws.Range("I84").Value = 3567
$ws.Range("J84").Value = 58299.875
$ws.Range("K84").Value = 35670
$ws.Range("L84").Value = 582998.75
$ws.Range("M84").Value = -30366
$ws.Range("N84").Value = -593606.75

$ws.Range("H100").Value = 7939106.5
$ws.Range("I100").Value = 17861494
$ws.Range("J100").Value = 1197.2
$ws.Range("K100").Value = 35722988
$ws.Range("L100").Value = 2394.4
$ws.Range("M100").Value = -35722447
$ws.Range("N100").Value = -3476.4

$ws.Range("H102").Value = 60500
$ws.Range("J102").Value = 60500
$ws.Range("L102").Value = 60500
$ws.Range("N102").Value = -66990

$ws.Range("H106").Value = 48500
$ws.Range("J106").Value = 48500
$ws.Range("L106").Value = 48500
$ws.Range("N106").Value = -51024

$ws.Range("H109").Value = 54641.75
$ws.Range("J109").Value = 54641.75
$ws.Range("L109").Value = 54641.75
$ws.Range("N109").Value = -57415.75

$ws.Range("H132").Value = 1319226.2
$ws.Range("I132").Value = 1470.6666
$ws.Range("J132").Value = 7249126.5
$ws.Range("K132").Value = 4411.9998
$ws.Range("L132").Value = 21747379.5
$ws.Range("M132").Value = -1881.9998
$ws.Range("N132").Value = -21752439.5
